$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1033.3158
$ws.Range("J17").Value = 1033.3158
$ws.Range("L17").Value = 3099.9474
$ws.Range("N17").Value = -3435.9474

$ws.Range("H76").Value = 4900
$ws.Range("I76").Value = 2350
$ws.Range("K76").Value = 2350
$ws.Range("M76").Value = -2035

$ws.Range("H79").Value = 4900
$ws.Range("I79").Value = 2350
$ws.Range("K79").Value = 2350
$ws.Range("M79").Value = -1258

$ws.Range("H86").Value = 202207.6
$ws.Range("I86").Value = 2758.5
$ws.Range("K86").Value = 2758.5
$ws.Range("M86").Value = -1635.5

$ws.Range("H89").Value = 202207.6
$ws.Range("I89").Value = 2758.5
$ws.Range("K89").Value = 13792.5
$ws.Range("M89").Value = -8176.5

$ws.Range("H138").Value = 5809
$ws.Range("I138").Value = 5374.5
$ws.Range("J138").Value = 5905.5557
$ws.Range("K138").Value = 16123.5
$ws.Range("L138").Value = 17716.6671
$ws.Range("M138").Value = -10983.5
$ws.Range("N138").Value = -27996.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1400
$ws.Range("I45").Value = 1400
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1400
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1023
$ws.Range("N45").ClearContents()

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H132").Value = 3346.4
$ws.Range("I132").Value = 3346.4
$ws.Range("K132").Value = 10039.2
$ws.Range("M132").Value = -7509.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1999.4
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

$ws.Range("H105").Value = 2600
$ws.Range("I105").Value = 2600
$ws.Range("J105").Value = 2600
$ws.Range("K105").Value = 2600
$ws.Range("L105").Value = 2600
$ws.Range("M105").Value = -853
$ws.Range("N105").Value = -6094

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2564.7778
$ws.Range("I132").Value = 2484.6667
$ws.Range("J132").Value = 2965.3333
$ws.Range("K132").Value = 7454.000100000001
$ws.Range("L132").Value = 8895.999899999999
$ws.Range("M132").Value = -4924.000100000001
$ws.Range("N132").Value = -13955.9999

$ws.Range("H141").Value = 37090.855
$ws.Range("J141").Value = 37090.855
$ws.Range("L141").Value = 37090.855
$ws.Range("N141").Value = -47450.855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 999
$ws.Range("J24").Value = 999
$ws.Range("L24").Value = 2997
$ws.Range("N24").Value = -3457

$ws.Range("H31").Value = 100
$ws.Range("J31").Value = 100
$ws.Range("L31").Value = 300
$ws.Range("N31").Value = -876

$ws.Range("H33").Value = 1098.5
$ws.Range("I33").Value = 1098
$ws.Range("K33").Value = 6588
$ws.Range("M33").Value = -6305

$ws.Range("H48").Value = 250
$ws.Range("J48").Value = 250
$ws.Range("L48").Value = 750
$ws.Range("N48").Value = -1250

$ws.Range("H99").Value = 1749.5
$ws.Range("I99").Value = 2500
$ws.Range("J99").Value = 999
$ws.Range("K99").Value = 7500
$ws.Range("L99").Value = 2997
$ws.Range("M99").Value = -5254
$ws.Range("N99").Value = -7489

$ws.Range("H106").Value = 18249.875
$ws.Range("I106").Value = 15666.667
$ws.Range("K106").Value = 47000.001
$ws.Range("M106").Value = -46054.001

$ws.Range("H107").Value = 527
$ws.Range("I107").Value = 582.6667
$ws.Range("J107").Value = 360
$ws.Range("K107").Value = 1748.0001
$ws.Range("L107").Value = 1080
$ws.Range("M107").Value = 171.9999
$ws.Range("N107").Value = -4920

$ws.Range("H112").Value = 11331.125
$ws.Range("I112").Value = 2662.25
$ws.Range("K112").Value = 7986.75
$ws.Range("M112").Value = -6878.75

$ws.Range("H131").Value = 2719.9
$ws.Range("J131").Value = 2911
$ws.Range("L131").Value = 8733
$ws.Range("N131").Value = -18813

$ws.Range("H132").Value = 3895.8
$ws.Range("J132").Value = 4500
$ws.Range("L132").Value = 40500
$ws.Range("N132").Value = -45560

$ws.Range("H137").Value = 3901.625
$ws.Range("J137").Value = 3901.625
$ws.Range("L137").Value = 11704.875
$ws.Range("N137").Value = -21904.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 26000
$ws.Range("J47").Value = 26000
$ws.Range("L47").Value = 26000
$ws.Range("N47").Value = -27136

$ws.Range("H70").Value = 5251.5
$ws.Range("I70").Value = 4999.5
$ws.Range("K70").Value = 4999.5
$ws.Range("M70").Value = -4729.5

$ws.Range("H73").Value = 5251.5
$ws.Range("I73").Value = 4999.5
$ws.Range("K73").Value = 4999.5
$ws.Range("M73").Value = -4063.5

$ws.Range("H102").Value = 2420.5
$ws.Range("I102").Value = 2585.375
$ws.Range("J102").Value = 1761
$ws.Range("K102").Value = 2585.375
$ws.Range("L102").Value = 1761
$ws.Range("M102").Value = -963.375
$ws.Range("N102").Value = -5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 545.2
$ws.Range("I16").Value = 545.2
$ws.Range("K16").Value = 545.2
$ws.Range("M16").Value = -375.2

$ws.Range("H136").Value = 3936.889
$ws.Range("I136").Value = 3936.889
$ws.Range("K136").Value = 11810.667
$ws.Range("M136").Value = -9260.667000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 19999
$ws.Range("I38").Value = 19999
$ws.Range("K38").Value = 19999
$ws.Range("M38").Value = -19526

$ws.Range("H42").Value = 15044
$ws.Range("I42").Value = 15044
$ws.Range("K42").Value = 15044
$ws.Range("M42").Value = -14666

$ws.Range("H107").Value = 508.75
$ws.Range("I107").Value = 216
$ws.Range("J107").Value = 801.5
$ws.Range("K107").Value = 648
$ws.Range("L107").Value = 2404.5
$ws.Range("M107").Value = 1272
$ws.Range("N107").Value = -6244.5

$ws.Range("H132").Value = 3125.5
$ws.Range("I132").Value = 3125.5
$ws.Range("K132").Value = 9376.5
$ws.Range("M132").Value = -6846.5
